$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.385.92'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.02%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.567.93'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.11%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.32'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.47%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3763'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.27%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.52'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.50%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3409'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07609'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.143'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.99%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.06'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.987'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.957'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.567.23'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.66%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001133'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.96'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.33%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.21%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.56'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.194'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.95'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.387.28'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.404'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.45%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.720'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -6.77%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.14'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.61%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '147.54'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.14%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.71%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.40'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.61%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.742.77'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.26%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.001'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.105'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.36%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '10.14'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.75%  '

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.441'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +10.00%  '

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'Stellar'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08516'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.59%  '

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02520'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.70%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2304'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06498'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.21%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.404'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.20%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.37'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6337'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.08%  '

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.02'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.38%  '

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.804'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.33%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5940'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.60%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'EOS'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.281'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.49%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.079'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.08%  '

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.34'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.53%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07311'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.51%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Flow'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.065'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.26%  '
